$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Odds Ratio / CI values for rows 2-20 (sensitivity: blood gas only)
$data = @{
    2  = @(0.9488411582465469, 0.407354605155226, 2.210112595240163)
    3  = @(1.253413971888949, 1.107642987878564, 1.418369097370815)
    4  = @(0.9723939511214528, 0.8381951595210121, 1.128078569098307)
    5  = @(0.9903870926872311, 0.9733128763468433, 1.007760831278606)
    6  = @(1.197627154790709, 1.184477668066314, 1.210922620629592)
    7  = @(0.9973503141322366, 0.9939883257945266, 1.000723673796239)
    8  = @(0.9917223146622363, 0.8557189263274573, 1.149341353965406)
    9  = @(6.57006563924765, 6.027333300987554, 7.161668410962123)
    10 = @(1.022393387768628, 0.8750852033297862, 1.194498816087377)
    11 = @(1.069636368392201, 0.91823621251535, 1.245999607718727)
    12 = @(0.9852563497576496, 0.9105167048164537, 1.066130988704323)
    13 = @(0.898572494923975, 0.7906746335718713, 1.021194426064135)
    14 = @(0.8838844303016808, 0.796355427133532, 0.9810339196680422)
    15 = @(0.9469882888737423, 0.8585160232803121, 1.044577847059251)
    16 = @(1.041734379825617, 0.8426191018382777, 1.287901634016061)
    17 = @(0.9343034340384859, 0.8201330176576371, 1.064367472180601)
    18 = @(0.9034324262804929, 0.7133777418229791, 1.144120570357789)
    19 = @(1.03448441586587, 0.9483984944080451, 1.128384337363697)
    20 = @(0.9810424703866519, 0.8335894966537385, 1.15457828171524)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
}

# Add new row 21 ("const")
$ws.Cells.Item(21, 1).Value = "const"
$ws.Cells.Item(20, 1).Copy()
$ws.Cells.Item(21, 1).PasteSpecial(-4122)

# B21 stays blank (matches the empty-string cell in the source data export)
$ws.Cells.Item(21, 2).Formula = "="""""

$ws.Cells.Item(21, 3).Value = 0.02458040599925118
$ws.Cells.Item(21, 4).Value = 0.01844754074416798
$ws.Cells.Item(21, 5).Value = 0.03275213577067364
